# The workbook's "document map" sheet listed a tutorial/document column
# for "御紹介@"棟梁".pptx (Introduction.pptx)" in column Q. This document
# reference was removed from the map (see issues #12 / #19), so the whole
# column is deleted - Excel shifts every column to its right (R:Y) one
# slot to the left (R->Q, S->R, ... X->W, Y->X), which is exactly what a
# real "delete entire column" action does.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q1").EntireColumn.Delete()
